$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 305, shifting existing rows 305:334 down to 306:335.
$ws.Rows.Item(305).Insert()

# Populate the newly inserted row 305 with the new weekly record.
$ws.Cells.Item(305, 1).Value = 10
$ws.Cells.Item(305, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(305, 3).Value = "La Araucanía"
$ws.Cells.Item(305, 4).Value = 44474
$ws.Cells.Item(305, 5).Value = 9
$ws.Cells.Item(305, 6).Value = "Fruta"
$ws.Cells.Item(305, 7).Value = 100101
$ws.Cells.Item(305, 8).Value = "Berries"
$ws.Cells.Item(305, 9).Value = 100101007
$ws.Cells.Item(305, 10).Value = "Kiwi"
$ws.Cells.Item(305, 11).Value = "Hayward"
$ws.Cells.Item(305, 12).Value = "Especial"
$ws.Cells.Item(305, 13).Value = 80
$ws.Cells.Item(305, 14).Value = 17000
$ws.Cells.Item(305, 15).Value = 17000
$ws.Cells.Item(305, 16).Value = 17000
$ws.Cells.Item(305, 17).Value = '$/caja 15 kilos'
$ws.Cells.Item(305, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(305, 19).Value = 1133
$ws.Cells.Item(305, 20).Value = 15
